$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the pictures anchored next to the item rows (10-12) that are about
# to be deleted, since picture anchors aren't auto-pruned by row deletion.
$ws.Shapes.Item("그림 1").Delete()
$ws.Shapes.Item("그림 15").Delete()
$ws.Shapes.Item("그림 3").Delete()

# Remove the hyperlink that lives on row 11 (item 5's link), since row
# deletion doesn't automatically prune hyperlinks anchored in removed rows.
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$F$11') {
        $hl.Delete()
    }
}

# Delete item rows 10-12 (the 4th, 5th, 6th material entries), along with
# their hyperlinks and row-specific style/font usage.
$ws.Rows("10:12").Delete()

# The footer note row (old row 14, now shifted up to row 11) should be
# cleared of its text and un-merged.
$ws.Range("A11:G11").UnMerge()
$ws.Range("A11:G11").ClearContents()

$ws.Range("D15").Select()
